# Remove console messages / display full error messages on console
# -> Replace the four placeholder "animal" quiz questions (cat/dog/tiger/lion)
#    in the Question sheet with the real quiz content, and move the active
#    selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Question")

# --- Row 2 : "cat" -> "DDD stands for?" -------------------------------
$ws.Cells.Item(2, 3).Value = "DDD stands for?"
$ws.Cells.Item(2, 5).Value = "DDD"
$ws.Cells.Item(2, 6).Value = "ACB"
$ws.Cells.Item(2, 7).Value = "HEF"
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 0

# --- Row 3 : "dog" -> "Who can have a alpha?" --------------------------
$ws.Cells.Item(3, 5).Value = "Manager"
$ws.Cells.Item(3, 6).Value = "Driver"
$ws.Cells.Item(3, 7).Value = "You"
$ws.Cells.Item(3, 8).Value = "SME"
$ws.Cells.Item(3, 3).Value = "Who can have a alpha?"
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0

# --- Row 4 : "tiger" -> "Which offshore team handle issue related a driver?"
$ws.Cells.Item(4, 5).Value = "Backoffice"
$ws.Cells.Item(4, 6).Value = "Order"
$ws.Cells.Item(4, 7).Value = "Vehicle"
$ws.Cells.Item(4, 8).Value = "Operations"
$ws.Cells.Item(4, 3).Value = "Which offshore team handle issue related a driver?"
$ws.Cells.Item(4, 9).Value = ""
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = ""

# --- Row 5 : "lion" -> "How are all SMEs for Vehicle?" ------------------
$ws.Cells.Item(5, 3).Value = "How are all SMEs for Vehicle?"
$ws.Cells.Item(5, 4).Value = "MCA"
$ws.Cells.Item(5, 5).Value = "Jamie"
$ws.Cells.Item(5, 6).Value = "Anitha"
$ws.Cells.Item(5, 7).Value = "Maria"
$ws.Cells.Item(5, 8).Value = "Vlad"
$ws.Cells.Item(5, 9).Value = "Prasanna"
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1
$ws.Cells.Item(5, 14).Value = 0

# Move the active selection (was C11)
$ws.Range("L6").Select()
